# kingdom finance rampung dan data dummy
# - Rework amd_users "access" row: new wording + selection change
# - Re-select full range on amd_articleCategory
# - Insert "image" column into amd_article table
# - Add three new reference tables: amd_partner, amd_social, amd_website

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. amd_users: update the "access" description text & selection
# ---------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("amd_users")
$wsUsers.Range("G2").Value = "0=all, 1=all except edit or create user, 2=spe"
$wsUsers.Range("G2").Select() | Out-Null

# ---------------------------------------------------------------
# 2. amd_articleCategory: select the whole data range
# ---------------------------------------------------------------
$wsArticleCategory = $wb.Worksheets.Item("amd_articleCategory")
$wsArticleCategory.Range("A1:N3").Select() | Out-Null

# ---------------------------------------------------------------
# 3. amd_article: insert a new "image" column before "priority"
# ---------------------------------------------------------------
$wsArticle = $wb.Worksheets.Item("amd_article")
$wsArticle.Columns.Item(7).Insert() | Out-Null

$wsArticle.Range("G1").Value = "image"
$wsArticle.Range("G2").Value = "nullable"
$wsArticle.Range("G3").Value = "text"

$wsArticle.Range("H1:H3").Select() | Out-Null

# ---------------------------------------------------------------
# 4. Add new sheet amd_partner
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPartner = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsPartner.Name = "amd_partner"

$wsPartner.Range("A1").Value = "id*"
$wsPartner.Range("B1").Value = "name"
$wsPartner.Range("C1").Value = "link"
$wsPartner.Range("D1").Value = "image"
$wsPartner.Range("E1").Value = "priority"
$wsPartner.Range("F1").Value = "actor"
$wsPartner.Range("G1").Value = "flag_publish"

$wsPartner.Range("A2").Value = "primary"
$wsPartner.Range("C2").Value = "nullable"
$wsPartner.Range("E2").Value = "default: 0"
$wsPartner.Range("F2").Value = "foreign user.id"
$wsPartner.Range("G2").Value = "1 = publish, 0=unpublish"

$wsPartner.Range("A3").Value = "increment"
$wsPartner.Range("B3").Value = "string"
$wsPartner.Range("C3").Value = "text"
$wsPartner.Range("D3").Value = "text"
$wsPartner.Range("E3").Value = "integer"
$wsPartner.Range("F3").Value = "integer"
$wsPartner.Range("G3").Value = "bool"

$wsPartner.Columns.Item(5).ColumnWidth = 14.140625
$wsPartner.Columns("F:G").ColumnWidth = 22.7109375

$wsPartner.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------
# 5. Add new sheet amd_social (same shape as amd_partner)
# ---------------------------------------------------------------
$wsSocial = $wb.Worksheets.Add([System.Type]::Missing, $wsPartner)
$wsSocial.Name = "amd_social"

$wsSocial.Range("A1").Value = "id*"
$wsSocial.Range("B1").Value = "name"
$wsSocial.Range("C1").Value = "link"
$wsSocial.Range("D1").Value = "image"
$wsSocial.Range("E1").Value = "priority"
$wsSocial.Range("F1").Value = "actor"
$wsSocial.Range("G1").Value = "flag_publish"

$wsSocial.Range("A2").Value = "primary"
$wsSocial.Range("C2").Value = "nullable"
$wsSocial.Range("E2").Value = "default: 0"
$wsSocial.Range("F2").Value = "foreign user.id"
$wsSocial.Range("G2").Value = "1 = publish, 0=unpublish"

$wsSocial.Range("A3").Value = "increment"
$wsSocial.Range("B3").Value = "string"
$wsSocial.Range("C3").Value = "text"
$wsSocial.Range("D3").Value = "text"
$wsSocial.Range("E3").Value = "integer"
$wsSocial.Range("F3").Value = "integer"
$wsSocial.Range("G3").Value = "bool"

$wsSocial.Range("M17").Select() | Out-Null

# ---------------------------------------------------------------
# 6. Add new sheet amd_website
# ---------------------------------------------------------------
$wsWebsite = $wb.Worksheets.Add([System.Type]::Missing, $wsSocial)
$wsWebsite.Name = "amd_website"

$wsWebsite.Range("A1").Value = "id*"
$wsWebsite.Range("B1").Value = "website"
$wsWebsite.Range("C1").Value = "company"
$wsWebsite.Range("D1").Value = "brief"
$wsWebsite.Range("E1").Value = "address"
$wsWebsite.Range("F1").Value = "phone"
$wsWebsite.Range("G1").Value = "google_map_url"
$wsWebsite.Range("H1").Value = "actor"
$wsWebsite.Range("I1").Value = "flag_publish"

$wsWebsite.Range("A2").Value = "primary"
$wsWebsite.Range("B2").Value = "nullable"
$wsWebsite.Range("C2").Value = "nullable"
$wsWebsite.Range("D2").Value = "nullable"
$wsWebsite.Range("E2").Value = "nullable"
$wsWebsite.Range("F2").Value = "nullable"
$wsWebsite.Range("G2").Value = "nullable"
$wsWebsite.Range("H2").Value = "foreign user.id"
$wsWebsite.Range("I2").Value = "1 = publish, 0=unpublish"

$wsWebsite.Range("A3").Value = "increment"
$wsWebsite.Range("B3").Value = "string"
$wsWebsite.Range("C3").Value = "string"
$wsWebsite.Range("D3").Value = "string"
$wsWebsite.Range("E3").Value = "string"
$wsWebsite.Range("F3").Value = "text"
$wsWebsite.Range("G3").Value = "text"
$wsWebsite.Range("H3").Value = "integer"
$wsWebsite.Range("I3").Value = "bool"

$wsWebsite.Columns.Item(7).ColumnWidth = 16.28515625

$wsWebsite.Range("J9").Select() | Out-Null

# ---------------------------------------------------------------
# 7. Final workbook view: amd_social active, amd_articleCategory first visible
# ---------------------------------------------------------------
$wsArticleCategory.Activate() | Out-Null
$wsSocial.Activate() | Out-Null
$wsSocial.Select() | Out-Null
